$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.608.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.506.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'169.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.57%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.618"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.504.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.05%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.88%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.50%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.579"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'47.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.54%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.89%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.086.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.15%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'8.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.72%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'613.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -8.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.501.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.613.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.92%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.97%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'17.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.88%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -9.41%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -2.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'15.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'95.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.77%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.27%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.42%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'33.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.59%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.89%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.35%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.94%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.45%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'554.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'10.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -3.26%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'56.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -4.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.05%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Kaspa"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.142"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.62%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'VeChain"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.0449"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.339.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.325"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.63%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'32.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.49%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0" + [char]0x2083 + "0698"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.16%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.01%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.32%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'135.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'5.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +7.45%  "
$ws.Range("E51").Style = "Normal"
